# update week 4 timesheet
# Fill in the two new tasks that were logged on row 11 and row 12
# (Task / Date / Start time / End time / computed duration).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: "SRS draft preperation" - 2017-08-20, 08:00 -> 12:00
$ws.Range("C11").Value = "SRS draft preperation"
$ws.Range("D11").Value = 42967
$ws.Range("E11").Value = 0.33333333333333331
$ws.Range("F11").Value = 0.5
$ws.Range("G11").Formula = "=(F11-E11)"

# Row 12: "Ant tool research for build file" - 2017-08-20, 16:00 -> 17:30
$ws.Range("C12").Value = "Ant tool research for build file"
$ws.Range("D12").Value = 42967
$ws.Range("E12").Value = 0.66666666666666663
$ws.Range("F12").Value = 0.72916666666666663
$ws.Range("G12").Formula = "=(F12-E12)"

# Move the selection/scroll position to where the author left off editing.
$ws.Range("E13").Select()
$excel.ActiveWindow.ScrollRow = 10
